$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 198602002.3250627
$ws.Range("G2").Value = 198602169.3799862

$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 9.226618575922256
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("G3").Value = 71530.92951483269
